# Apply odds updates for 2025-06-11 FlashScore workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Cells.Item(6, 10).Value = 1.08    # J6  Odd_Over05_FT
$ws.Cells.Item(6, 12).Value = 1.36    # L6  Odd_Over15_FT

# Row 11
$ws.Cells.Item(11, 8).Value = 3.5     # H11 Odd_D_FT
$ws.Cells.Item(11, 9).Value = 3.85    # I11 Odd_A_FT
$ws.Cells.Item(11, 12).Value = 1.26   # L11 Odd_Over15_FT
$ws.Cells.Item(11, 13).Value = 3.15   # M11 Odd_Under15_FT
$ws.Cells.Item(11, 14).Value = 1.78   # N11 Odd_Over25_FT
$ws.Cells.Item(11, 15).Value = 1.83   # O11 Odd_Under25_FT
$ws.Cells.Item(11, 18).Value = 1.7    # R11 Odd_BTTS_Yes
$ws.Cells.Item(11, 19).Value = 1.93   # S11 Odd_BTTS_No
$ws.Cells.Item(11, 20).Value = 7.6    # T11 Odd_CS_1-0
$ws.Cells.Item(11, 21).Value = 9      # U11 Odd_CS_2-0
$ws.Cells.Item(11, 24).Value = 14     # X11 Odd_CS_3-1
$ws.Cells.Item(11, 25).Value = 24     # Y11 Odd_CS_3-2
$ws.Cells.Item(11, 26).Value = 10.5   # Z11 Odd_CS_0-0
$ws.Cells.Item(11, 27).Value = 6.8    # AA11 Odd_CS_1-1
$ws.Cells.Item(11, 28).Value = 14     # AB11 Odd_CS_2-2
$ws.Cells.Item(11, 29).Value = 60     # AC11 Odd_CS_3-3
$ws.Cells.Item(11, 30).Value = 450    # AD11 Odd_CS_4-4
$ws.Cells.Item(11, 31).Value = 11.5   # AE11 Odd_CS_0-1
$ws.Cells.Item(11, 32).Value = 22     # AF11 Odd_CS_0-2
$ws.Cells.Item(11, 33).Value = 13     # AG11 Odd_CS_1-2
$ws.Cells.Item(11, 35).Value = 35     # AI11 Odd_CS_1-3
$ws.Cells.Item(11, 36).Value = 40     # AJ11 Odd_CS_2-3

# Row 12
$ws.Cells.Item(12, 14).Value = 1.45   # N12 Odd_Over25_FT
$ws.Cells.Item(12, 15).Value = 2.37   # O12 Odd_Under25_FT
$ws.Cells.Item(12, 18).Value = 1.52   # R12 Odd_BTTS_Yes
$ws.Cells.Item(12, 20).Value = 10.25  # T12 Odd_CS_1-0
$ws.Cells.Item(12, 21).Value = 9.75   # U12 Odd_CS_2-0
$ws.Cells.Item(12, 26).Value = 17.5   # Z12 Odd_CS_0-0
$ws.Cells.Item(12, 31).Value = 18     # AE12 Odd_CS_0-1

# Row 13
$ws.Cells.Item(13, 7).Value = 2.45    # G13 Odd_H_FT
$ws.Cells.Item(13, 9).Value = 2.75    # I13 Odd_A_FT

# Row 14
$ws.Cells.Item(14, 7).Value = 1.62    # G14 Odd_H_FT
$ws.Cells.Item(14, 9).Value = 4.5     # I14 Odd_A_FT
$ws.Cells.Item(14, 18).Value = 1.72   # R14 Odd_BTTS_Yes
$ws.Cells.Item(14, 19).Value = 2      # S14 Odd_BTTS_No
$ws.Cells.Item(14, 21).Value = 8.25   # U14 Odd_CS_2-0
$ws.Cells.Item(14, 23).Value = 12.5   # W14 Odd_CS_3-0
$ws.Cells.Item(14, 24).Value = 12.5   # X14 Odd_CS_3-1
$ws.Cells.Item(14, 27).Value = 7.7    # AA14 Odd_CS_1-1
$ws.Cells.Item(14, 28).Value = 15     # AB14 Odd_CS_2-2
$ws.Cells.Item(14, 29).Value = 60     # AC14 Odd_CS_3-3
$ws.Cells.Item(14, 32).Value = 28     # AF14 Odd_CS_0-2
$ws.Cells.Item(14, 33).Value = 15     # AG14 Odd_CS_1-2
$ws.Cells.Item(14, 35).Value = 40     # AI14 Odd_CS_1-3
$ws.Cells.Item(14, 36).Value = 40     # AJ14 Odd_CS_2-3
